# Apply the edit described by the commit: add "savee" dataset row to the
# samples_retained sheet, and switch the active sheet/selection from
# "neutral" (row 24 area) back to "samples_retained".

$wb = $excel.ActiveWorkbook

$wsSamples = $wb.Worksheets.Item("samples_retained")
$wsNeutral = $wb.Worksheets.Item("neutral")

# --- Fill in the new "savee" row (row 24) on samples_retained ---
$wsSamples.Range("B24").Value = "acted"
$wsSamples.Range("C24").Value = 61
$wsSamples.Range("D24").Value = 306
$wsSamples.Range("E24").Value = 121
$wsSamples.Range("F24").Value = "English"
$wsSamples.Range("H24").Formula = "=4+7"
$wsSamples.Range("I24").Value = "anger, disgust, fear, happiness, neutral, sadness, surprise"
$wsSamples.Range("J24").Value = "added more unique speakers from MetaData folder; 4 main male speakers otherwise"

# --- Move the active sheet / selection back to samples_retained ---
$wsNeutral.Range("A5").Select()
$wsSamples.Activate()
$wsSamples.Range("F24").Select()
